# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Sat Oct 12 09:42:55 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "62.742.07" / "576.01". Excel's COM layer
# auto-coerces plain decimal-looking strings (e.g. "576.01") into Numbers,
# which would corrupt both the stored type and the text (float rounding).
# Force text storage via a temporary "@" (Text) format, then restore the
# default "Normal" style so the cell ends up styled exactly like its
# untouched neighbours (no leftover per-cell style index).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.742.07"
$ws.Range("E2").Value = "  +3.13%  "
Set-TextValue "D3" "2.445.50"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue "D5" "576.01"
$ws.Range("E5").Value = "  +1.79%  "
Set-TextValue "D6" "145.63"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue "D8" "0.540"
$ws.Range("E8").Value = "  +0.12%  "
Set-TextValue "D9" "2.444.34"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("E13").Value = "  +2.18%  "
Set-TextValue "D14" "28.38"
$ws.Range("E14").Value = "  +7.48%  "
$ws.Range("E15").Value = "  +5.40%  "
Set-TextValue "D16" "2.890.33"
$ws.Range("E16").Value = "  +1.87%  "
Set-TextValue "D17" "62.651.23"
$ws.Range("E17").Value = "  +3.43%  "
Set-TextValue "D18" "2.444.59"
$ws.Range("E18").Value = "  +1.44%  "
Set-TextValue "D19" "7.93"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("E20").Value = "  +3.00%  "
Set-TextValue "D21" "330.13"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("E23").Value = "  +8.00%  "
$ws.Range("E24").Value = "  +0.09%  "
Set-TextValue "D25" "66.29"
$ws.Range("E25").Value = "  +1.85%  "
Set-TextValue "D26" "646.83"
$ws.Range("E26").Value = "  +10.83%  "
$ws.Range("E27").Value = "  +17.80%  "
Set-TextValue "D28" "8.53"
$ws.Range("E28").Value = "  +3.79%  "
Set-TextValue "D29" "0.0₃0988"
$ws.Range("E29").Value = "  +5.12%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D30" "2.567.29"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("B31").Value = "BabyDogeCoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D31" "0.0₆0488"
$ws.Range("E31").Value = "  +74.92%  "
Set-TextValue "D32" "8.19"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("E33").Value = "  +6.68%  "
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +2.97%  "
$ws.Range("E39").Value = "  +6.12%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D40" "153.39"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D41" "0.374"
$ws.Range("E41").Value = "  +0.70%  "
Set-TextValue "D42" "18.74"
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("E43").Value = "  +8.12%  "
$ws.Range("E44").Value = "  +4.15%  "
Set-TextValue "D45" "42.46"
$ws.Range("E46").Value = "  +0.02%  "
Set-TextValue "D47" "14.94"
$ws.Range("E47").Value = "  +27.45%  "
Set-TextValue "D48" "145.09"
$ws.Range("E48").Value = "  +2.55%  "
Set-TextValue "D49" "3.63"
$ws.Range("E49").Value = "  +3.40%  "
$ws.Range("E50").Value = "  +5.79%  "
$ws.Range("E51").Value = "  +2.38%  "
